$wb = $excel.ActiveWorkbook

$wsRecon = $wb.Worksheets.Item("Reconciliation Peer Review")
$wsSpec  = $wb.Worksheets.Item("Specification QA")

# --- Specification QA sheet: update the scrolled/frozen view & selection ---
# (scroll the frozen pane down so row 19 is the first visible row, and move
# the selected cell from G24 to A24, per the recorded view state)
$wsSpec.Activate()
$wsSpec.Application.Goto($wsSpec.Range("C19"), $true)
$wsSpec.Range("A24").Select()

# --- Reconciliation Peer Review sheet: add the new disposition row ---
$wsRecon.Range("A2").Value = 207
$wsRecon.Range("B2").Value = "The ballot comment was misunderstood.  Conformance still needs to declare what search parameters will be supported.  However, Profile should be used to *define* search parameters - you want the search parameters to be able to be defined the same place you define the resource element or extension being searched on."
$wsRecon.Range("C2").Value = "Y"
$wsRecon.Range("D2").Value = "for San Antonio"

# Row 2 wraps onto multiple lines in Excel (ht="90") once the comment is entered
$wsRecon.Rows.Item(2).RowHeight = 90

# This becomes the active / selected sheet and range (tab moves from the 3rd
# to the 2nd sheet, i.e. activeTab 2 -> 1)
$wsRecon.Activate()
$wsRecon.Range("A2:D2").Select()
